$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# ---------------------------------------------------------------------------
# Update the "Data" (GDP per Capita) values for the existing years 1950-2008
# (rows 2-60). Values are numeric-looking text that must remain stored as
# text (shared-string) cells, matching the original file's cell typing.
# Formatting the cell as Text ("@") before assigning forces Excel to keep
# the value as a string instead of auto-converting it to a number; clearing
# the format afterwards removes the temporary Text number-format again so
# the cell is left with the workbook's default (unstyled) formatting.
# ---------------------------------------------------------------------------
$newValues = @("1693", "1726", "1758", "1790", "1822", "1854", "1884", "1914", "1946", "1977", "1999", "2020", "2039", "2059", "2079", "2099", "2157", "2216", "2275", "2335", "2396", "2546", "2582", "2950", "2848", "2745", "2794", "2517", "2595", "2807", "3075", "3625", "4355", "4473", "4664", "4484", "4058", "3953", "3908", "3865", "3794", "3862.57433398122", "3940.97172910264", "3874.16978461464", "3635.98917024344", "3763.59580118657", "3988.87730463448", "3899.76268542608", "4042.77828134958", "3925.68643829626", "4159.22273739588", "4256.65537117082", "4458.63422368835", "4467.9039644057", "4560.51251418364", "4820.70728800517", "5016.05763935475", "4844.9702396029", "5081.13158781301")

for ($i = 0; $i -lt $newValues.Length; $i++) {
    $row = $i + 2
    $cell = $ws.Cells.Item($row, 5)
    $cell.NumberFormat = "@"
    $cell.Value = $newValues[$i]
    $cell.ClearFormats()
}

# ---------------------------------------------------------------------------
# Append the new years 2009-2016 (rows 61-68) with their GDP per Capita data.
# ---------------------------------------------------------------------------
$newYears = @(2009, 2010, 2011, 2012, 2013, 2014, 2015, 2016)
$newRowValues = @("5379.00417690191", "5738.16299309802", "5875", "5949", "6031", "6321", "6272", "5974")

for ($i = 0; $i -lt $newYears.Length; $i++) {
    $row = 61 + $i
    $ws.Cells.Item($row, 1).Value = 178
    $ws.Cells.Item($row, 2).Value = "Congo"
    $ws.Cells.Item($row, 3).Value = "GDP per Capita"
    $ws.Cells.Item($row, 4).Value = $newYears[$i]

    $cell = $ws.Cells.Item($row, 5)
    $cell.NumberFormat = "@"
    $cell.Value = $newRowValues[$i]
    $cell.ClearFormats()
}
